# Updated symbol list on Sun Jan 15 05:32:04 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto rows on Sheet1. The source cells are stored as plain text (e.g.
# "298.58", "-2.35%"), so we force each target cell to Text format before
# writing the new value - this keeps Excel from re-interpreting the string
# as a number/percentage and losing the exact formatting (trailing zeros,
# sign, "%" suffix, etc.) that the original data uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2';  Value = '298.58' },
    @{ Cell = 'E2';  Value = '-2.35%' },
    @{ Cell = 'D3';  Value = '31.87' },
    @{ Cell = 'E3';  Value = '-1.39%' },
    @{ Cell = 'D4';  Value = '5.094' },
    @{ Cell = 'E4';  Value = '-4.54%' },
    @{ Cell = 'D5';  Value = '0.07529' },
    @{ Cell = 'E5';  Value = '1.81%' },
    @{ Cell = 'D6';  Value = '7.765' },
    @{ Cell = 'E6';  Value = '0.10%' },
    @{ Cell = 'D7';  Value = '1.739' },
    @{ Cell = 'E7';  Value = '12.20%' },
    @{ Cell = 'D8';  Value = '3.790' },
    @{ Cell = 'E8';  Value = '2.47%' },
    @{ Cell = 'D9';  Value = '0.9280' },
    @{ Cell = 'E9';  Value = '1.97%' },
    @{ Cell = 'E10'; Value = '2.17%' },
    @{ Cell = 'D11'; Value = '0.07369' },
    @{ Cell = 'E11'; Value = '-2.13%' },
    @{ Cell = 'D12'; Value = '0.07971' },
    @{ Cell = 'E12'; Value = '-0.11%' },
    @{ Cell = 'D13'; Value = '0.03054' },
    @{ Cell = 'E13'; Value = '0.90%' },
    @{ Cell = 'D14'; Value = '0.09888' },
    @{ Cell = 'E14'; Value = '0.26%' },
    @{ Cell = 'D15'; Value = '0.001495' },
    @{ Cell = 'E15'; Value = '-1.68%' },
    @{ Cell = 'D16'; Value = '0.04646' },
    @{ Cell = 'E16'; Value = '1.84%' },
    @{ Cell = 'D17'; Value = '0.006546' },
    @{ Cell = 'E17'; Value = '2.88%' },
    @{ Cell = 'E18'; Value = '-0.74%' },
    @{ Cell = 'D19'; Value = '2.219' },
    @{ Cell = 'E19'; Value = '-0.92%' },
    @{ Cell = 'D21'; Value = '0.1316' },
    @{ Cell = 'E21'; Value = '-0.95%' },
    @{ Cell = 'D22'; Value = '4.558' },
    @{ Cell = 'E22'; Value = '8.07%' },
    @{ Cell = 'D23'; Value = '0.1548' },
    @{ Cell = 'E23'; Value = '-5.08%' },
    @{ Cell = 'E24'; Value = '-0.08%' },
    @{ Cell = 'D25'; Value = '0.004421' },
    @{ Cell = 'E25'; Value = '-1.78%' },
    @{ Cell = 'D26'; Value = '0.0001399' },
    @{ Cell = 'E26'; Value = '19.51%' },
    @{ Cell = 'D27'; Value = '0.0001934' },
    @{ Cell = 'E27'; Value = '7.25%' },
    @{ Cell = 'E39'; Value = '1.78%' },
    @{ Cell = 'D40'; Value = '0.04546' },
    @{ Cell = 'E40'; Value = '0.90%' },
    @{ Cell = 'D41'; Value = '0.007045' },
    @{ Cell = 'E41'; Value = '-5.21%' },
    @{ Cell = 'E42'; Value = '-2.63%' },
    @{ Cell = 'D43'; Value = '0.002058' },
    @{ Cell = 'E43'; Value = '-8.96%' },
    @{ Cell = 'E44'; Value = '-7.65%' },
    @{ Cell = 'D45'; Value = '0.00006050' },
    @{ Cell = 'E45'; Value = '-1.68%' },
    @{ Cell = 'D46'; Value = '0.7116' },
    @{ Cell = 'E46'; Value = '-62.40%' },
    @{ Cell = 'D47'; Value = '0.01297' },
    @{ Cell = 'E47'; Value = '-0.40%' }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = '@'
    $range.Value = $update.Value
}

$wb.Save()
